{"js": "// The article used to open with a bookmarked \"Heading1\" title paragraph\n// (\"Peter Maurin Farm\") followed by a bold \"By Dorothy Day\" byline\n// paragraph. Articles are now downloaded with pandoc-style title blocks:\n// a \"Title\" styled paragraph holding just the title, followed by an\n// \"Authors\" styled paragraph holding just the author name(s) -- each\n// word (and the spaces between them) kept as its own run, matching the\n// pandoc docx writer's output.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\n// First two paragraphs of the body: the \"Peter Maurin Farm\" heading and\n// the \"By Dorothy Day\" byline right after it.\nconst titlePara = body.paragraphs.items[0];\nconst bylinePara = body.paragraphs.items[1];\ntitlePara.load(\"text\");\nbylinePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text.trim() !== \"Peter Maurin Farm\" || bylinePara.text.trim() !== \"By Dorothy Day\") {\n  throw new Error(\"Unexpected document structure; aborting to avoid corrupting the wrong content.\");\n}\n\n// Grab a range spanning from the very start of the title paragraph to the\n// very end of the byline paragraph so a single replace swaps both\n// paragraphs for the new title-block paragraphs in one shot.\nconst startRange = titlePara.getRange(Word.RangeLocation.start);\nconst endRange = bylinePara.getRange(Word.RangeLocation.end);\nconst replaceRange = startRange.expandTo(endRange);\nawait context.sync();\n\nconst newOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n  <w:pPr><w:pStyle w:val=\"Title\"/></w:pPr>\n  <w:r><w:t xml:space=\"preserve\">Peter</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\">Maurin</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\">Farm</w:t></w:r>\n</w:p>\n<w:p>\n  <w:pPr><w:pStyle w:val=\"Authors\"/></w:pPr>\n  <w:r><w:t xml:space=\"preserve\">Dorothy</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\">Day</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\nreplaceRange.insertOoxml(newOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The article used to open with a bookmarked \"Heading1\" title paragraph\n# (\"Peter Maurin Farm\") followed by a bold \"By Dorothy Day\" byline\n# paragraph. Articles are now downloaded with pandoc-style title blocks:\n# a \"Title\" styled paragraph holding just the title, followed by an\n# \"Authors\" styled paragraph holding just the author name(s) -- each\n# word (and the spaces between them) kept as its own run, matching the\n# pandoc docx writer's output.\n\n$d = $word.ActiveDocument\n\n# First two paragraphs of the body: the \"Peter Maurin Farm\" heading and\n# the \"By Dorothy Day\" byline right after it.\n$titlePara = $d.Paragraphs.Item(1)\n$bylinePara = $d.Paragraphs.Item(2)\n\nif ($titlePara.Range.Text.Trim() -ne \"Peter Maurin Farm\" -or $bylinePara.Range.Text.Trim() -ne \"By Dorothy Day\") {\n    throw \"Unexpected document structure; aborting to avoid corrupting the wrong content.\"\n}\n\n# Range spanning from the very start of the title paragraph to the very\n# end of the byline paragraph, so a single XML insert swaps both\n# paragraphs for the new title-block paragraphs in one shot.\n$replaceRange = $d.Range($titlePara.Range.Start, $bylinePara.Range.End)\n\n$newOoxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n  <w:pPr><w:pStyle w:val=\"Title\"/></w:pPr>\n  <w:r><w:t xml:space=\"preserve\">Peter</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\">Maurin</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\">Farm</w:t></w:r>\n</w:p>\n<w:p>\n  <w:pPr><w:pStyle w:val=\"Authors\"/></w:pPr>\n  <w:r><w:t xml:space=\"preserve\">Dorothy</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\">Day</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$replaceRange.InsertXML($newOoxml)\n"}
